$wb = $excel.ActiveWorkbook

# --- Jornada sheet: clear the stray empty/styled cell B4, move selection ---
$wsJornada = $wb.Worksheets.Item("Jornada")
$wsJornada.Range("B4").Clear()

# --- Maquinas sheet: add prioritized printing process rows, fix Pegado machine name ---
$wsMaquinas = $wb.Worksheets.Item("Maquinas")
$wsMaquinas.Range("B5").Value = "Pegadora 1"
$wsMaquinas.Range("A7").Value = "Impresión Offset"
$wsMaquinas.Range("A8").Value = "Impresión Flexo"

# --- Update selections on each sheet ---
$wsJornada.Activate()
[void]$wsJornada.Range("E16").Select()

$wsMaquinas.Activate()
[void]$wsMaquinas.Range("C7").Select()
